# Update odds values for Jogos da Semana FlashScore 2025-05-30
# (values refreshed per latest FlashScore snapshot)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4
$ws.Range("G4").Value = 24
$ws.Range("H4").Value = 7.4
$ws.Range("I4").Value = 1.09
$ws.Range("N4").Value = 1.3
$ws.Range("O4").Value = 2.92
$ws.Range("R4").Value = 2.18
$ws.Range("S4").Value = 1.53
$ws.Range("T4").Value = 90
$ws.Range("U4").Value = 400
$ws.Range("Y4").Value = 250
$ws.Range("Z4").Value = 20
$ws.Range("AA4").Value = 19
$ws.Range("AC4").Value = 175
$ws.Range("AD4").Value = 10.25
$ws.Range("AE4").Value = 6.7
$ws.Range("AF4").Value = 12
$ws.Range("AG4").Value = 6.3
$ws.Range("AH4").Value = 11
$ws.Range("AI4").Value = 37

# Row 5
$ws.Range("G5").Value = 6.7
$ws.Range("H5").Value = 4.45
$ws.Range("I5").Value = 1.4
$ws.Range("L5").Value = 1.21
$ws.Range("M5").Value = 3.55
$ws.Range("N5").Value = 1.55
$ws.Range("O5").Value = 2.15
$ws.Range("R5").Value = 1.75
$ws.Range("S5").Value = 1.85
$ws.Range("T5").Value = 21
$ws.Range("U5").Value = 50
$ws.Range("V5").Value = 21
$ws.Range("W5").Value = 150
$ws.Range("Y5").Value = 55
$ws.Range("Z5").Value = 14.5
$ws.Range("AA5").Value = 9
$ws.Range("AB5").Value = 17.5
$ws.Range("AC5").Value = 70
$ws.Range("AD5").Value = 7.9
$ws.Range("AE5").Value = 7.2
$ws.Range("AH5").Value = 10.75
$ws.Range("AI5").Value = 23
$ws.Range("AJ5").Value = 500

# Row 6
$ws.Range("G6").Value = 1.2
$ws.Range("H6").Value = 5.8
$ws.Range("I6").Value = 12
$ws.Range("N6").Value = 1.44
$ws.Range("O6").Value = 2.42
$ws.Range("R6").Value = 1.98
$ws.Range("S6").Value = 1.65
$ws.Range("T6").Value = 8.5
$ws.Range("U6").Value = 6.4
$ws.Range("W6").Value = 7.1
$ws.Range("X6").Value = 10.5
$ws.Range("Y6").Value = 29
$ws.Range("Z6").Value = 16.5
$ws.Range("AA6").Value = 12.5
$ws.Range("AB6").Value = 26
$ws.Range("AC6").Value = 120
$ws.Range("AD6").Value = 35
$ws.Range("AE6").Value = 110
$ws.Range("AF6").Value = 40
$ws.Range("AG6").Value = 450
$ws.Range("AI6").Value = 120
$ws.Range("AJ6").Value = 900

# Row 7
$ws.Range("G7").Value = 1.17
$ws.Range("H7").Value = 6.1
$ws.Range("I7").Value = 16
$ws.Range("L7").Value = 1.24
$ws.Range("M7").Value = 3.65
$ws.Range("N7").Value = 1.72
$ws.Range("P7").Value = 1.35
$ws.Range("Q7").Value = 2.92
$ws.Range("S7").Value = 1.36
$ws.Range("U7").Value = 4.8
$ws.Range("W7").Value = 5.7
$ws.Range("AA7").Value = 13.5
$ws.Range("AC7").Value = 400
$ws.Range("AD7").Value = 30
$ws.Range("AE7").Value = 150
$ws.Range("AF7").Value = 60
$ws.Range("AG7").Value = 900
$ws.Range("AH7").Value = 400
$ws.Range("AI7").Value = 300

# Row 8
$ws.Range("J8").Value = 1.05
$ws.Range("K8").Value = 8.75
$ws.Range("T8").Value = 11.75
$ws.Range("Z8").Value = 8.75
$ws.Range("AC8").Value = 55
$ws.Range("AG8").Value = 18.5

# Row 11
$ws.Range("H11").Value = 3.5
$ws.Range("N11").Value = 1.95
$ws.Range("O11").Value = 1.85
$ws.Range("T11").Value = 7
$ws.Range("Y11").Value = 29
$ws.Range("Z11").Value = 10
$ws.Range("AA11").Value = 7
$ws.Range("AH11").Value = 34
$ws.Range("AJ11").Value = 500

# Row 19
$ws.Range("H19").Value = 3.7
$ws.Range("T19").Value = 7.7
$ws.Range("AD19").Value = 13
$ws.Range("AE19").Value = 25
$ws.Range("AI19").Value = 45

# Row 23
$ws.Range("AD23").Value = 6.8

# Row 25
$ws.Range("G25").Value = 1.27
$ws.Range("H25").Value = 4.85
$ws.Range("I25").Value = 8.5
$ws.Range("N25").Value = 1.62
$ws.Range("O25").Value = 2.02
$ws.Range("T25").Value = 5.8
$ws.Range("U25").Value = 5.1
$ws.Range("V25").Value = 7.6
$ws.Range("W25").Value = 6.3
$ws.Range("X25").Value = 9.5
$ws.Range("Y25").Value = 25
$ws.Range("Z25").Value = 12
$ws.Range("AA25").Value = 8.5
$ws.Range("AD25").Value = 17.5
$ws.Range("AE25").Value = 45
$ws.Range("AF25").Value = 23
$ws.Range("AG25").Value = 175

# Row 26
$ws.Range("L26").Value = 1.29
$ws.Range("M26").Value = 3.5
$ws.Range("N26").Value = 1.95
$ws.Range("O26").Value = 1.85

# Row 31
$ws.Range("G31").Value = 2.25
$ws.Range("H31").Value = 3.3
$ws.Range("I31").Value = 3.1
$ws.Range("K31").Value = 12
$ws.Range("L31").Value = 1.25
$ws.Range("M31").Value = 3.75
$ws.Range("N31").Value = 1.8
$ws.Range("O31").Value = 2
$ws.Range("R31").Value = 1.67
$ws.Range("S31").Value = 2.1
$ws.Range("T31").Value = 9
$ws.Range("W31").Value = 21
$ws.Range("Z31").Value = 12
$ws.Range("AG31").Value = 34
$ws.Range("AH31").Value = 23
$ws.Range("AI31").Value = 29

# Row 37
$ws.Range("G37").Value = 2.27
$ws.Range("H37").Value = 3.2
$ws.Range("I37").Value = 2.95
$ws.Range("L37").Value = 1.31
$ws.Range("M37").Value = 2.87
$ws.Range("N37").Value = 1.93
$ws.Range("O37").Value = 1.7
$ws.Range("P37").Value = 1.4
$ws.Range("Q37").Value = 2.52
$ws.Range("R37").Value = 1.72
$ws.Range("S37").Value = 1.88
$ws.Range("T37").Value = 7.4
$ws.Range("U37").Value = 11
$ws.Range("V37").Value = 9
$ws.Range("W37").Value = 23
$ws.Range("X37").Value = 19.5
$ws.Range("Y37").Value = 30
$ws.Range("Z37").Value = 9
$ws.Range("AA37").Value = 6.2
$ws.Range("AB37").Value = 14
$ws.Range("AC37").Value = 65
$ws.Range("AD37").Value = 9
$ws.Range("AE37").Value = 15.5
$ws.Range("AF37").Value = 10.5
$ws.Range("AG37").Value = 37
$ws.Range("AH37").Value = 26
$ws.Range("AI37").Value = 35
$ws.Range("AJ37").Value = 500
